$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45171 -> 45172) for every data row (rows 2 through 358).
$lastRow = 358
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45172
